$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.567.09"
$ws.Range("E2").Value = "  +0.83%  "
$ws.Range("D3").Value = "'3.169.17"
$ws.Range("E3").Value = "  -0.55%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'572.09"
$ws.Range("E5").Value = "  +0.12%  "
$ws.Range("D6").Value = "'164.19"
$ws.Range("E6").Value = "  -2.74%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "'0.587"
$ws.Range("E8").Value = "  -4.09%  "
$ws.Range("E9").Value = "  -3.00%  "
$ws.Range("D10").Value = "'6.63"
$ws.Range("E10").Value = "  -1.21%  "
$ws.Range("D11").Value = "'0.384"
$ws.Range("E11").Value = "  -0.90%  "
$ws.Range("D12").Value = "'3.719.31"
$ws.Range("E12").Value = "  -0.44%  "
$ws.Range("E13").Value = "  -1.02%  "
$ws.Range("D14").Value = "'64.579.17"
$ws.Range("E14").Value = "  +0.71%  "
$ws.Range("D15").Value = "'25.33"
$ws.Range("E15").Value = "  -0.59%  "
$ws.Range("D16").Value = "'3.172.83"
$ws.Range("E16").Value = "  -0.39%  "
$ws.Range("D17").Value = "'0.0000155"
$ws.Range("E17").Value = "  -2.78%  "
$ws.Range("D18").Value = "'408.63"
$ws.Range("E18").Value = "  -1.88%  "
$ws.Range("D19").Value = "'12.76"
$ws.Range("E19").Value = "  -0.59%  "
$ws.Range("D20").Value = "'5.27"
$ws.Range("E20").Value = "  -1.86%  "
$ws.Range("E21").Value = "  -0.73%  "
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("D23").Value = "'68.93"
$ws.Range("D24").Value = "'0.485"
$ws.Range("E24").Value = "  -1.73%  "
$ws.Range("E25").Value = "  -2.65%  "
$ws.Range("D26").Value = "'0.0000102"
$ws.Range("E26").Value = "  -6.53%  "
$ws.Range("D27").Value = "'8.87"
$ws.Range("E27").Value = "  +0.73%  "
$ws.Range("E28").Value = "  +0.19%  "
$ws.Range("D29").Value = "'1.82"
$ws.Range("E29").Value = "  -1.34%  "
$ws.Range("D30").Value = "'21.24"
$ws.Range("E30").Value = "  -3.30%  "
$ws.Range("B31").Value = "Aptos"
$ws.Range("C31").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D31").Value = "'6.35"
$ws.Range("E31").Value = "  -0.83%  "
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").Value = "'4.90"
$ws.Range("E32").Value = "  -2.48%  "
$ws.Range("E33").Value = "  -1.31%  "
$ws.Range("D34").Value = "'156.34"
$ws.Range("E34").Value = "  +0.27%  "
$ws.Range("E35").Value = "  -2.21%  "
$ws.Range("E36").Value = "  -0.78%  "
$ws.Range("D37").Value = "'2.683.74"
$ws.Range("E37").Value = "  -2.64%  "
$ws.Range("E38").Value = "  -4.07%  "
$ws.Range("D39").Value = "'4.10"
$ws.Range("E39").Value = "  -2.32%  "
$ws.Range("E40").Value = "  -2.88%  "
$ws.Range("D41").Value = "'0.0620"
$ws.Range("E41").Value = "  -1.40%  "
$ws.Range("D42").Value = "'5.44"
$ws.Range("E42").Value = "  -4.80%  "
$ws.Range("D43").Value = "'0.0258"
$ws.Range("E43").Value = "  -2.09%  "
$ws.Range("D44").Value = "'292.09"
$ws.Range("E44").Value = "  -2.23%  "
$ws.Range("D45").Value = "'21.43"
$ws.Range("E45").Value = "  -3.21%  "
$ws.Range("D46").Value = "'1.00"
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("D47").Value = "'0.0987"
$ws.Range("E47").Value = "  -0.91%  "
$ws.Range("D48").Value = "'1.92"
$ws.Range("E48").Value = "  -8.67%  "
$ws.Range("D49").Value = "'10.46"
$ws.Range("E49").Value = "  +0.20%  "
$ws.Range("D50").Value = "'5.72"
$ws.Range("E50").Value = "  -1.65%  "
$ws.Range("D51").Value = "'0.878"
$ws.Range("E51").Value = "  -6.38%  "
